# Update 2021 conference championship target depth data for the Chiefs.
# "OFF" sheet = Chiefs offense target depth splits (Home/Road)
# "DEF" sheet = Chiefs defense target depth splits (Home/Road)

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsDef = $wb.Worksheets.Item("DEF")

# OFF sheet, row 2 ("H" row): B2:G2
$wsOff.Range("B2").Value = 691
$wsOff.Range("C2").Value = 482
$wsOff.Range("D2").Value = 159
$wsOff.Range("E2").Value = 77
$wsOff.Range("F2").Value = 8
$wsOff.Range("G2").Value = 4

# DEF sheet, row 2 ("H" row): B2:G2 (F2 is unchanged at 13)
$wsDef.Range("B2").Value = 639
$wsDef.Range("C2").Value = 443
$wsDef.Range("D2").Value = 144
$wsDef.Range("E2").Value = 53
$wsDef.Range("G2").Value = 10
